$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.618.74'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '2.233.95'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '271.25'
$ws.Range('E5').Value = '  +4.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '90.78'
$ws.Range('E6').Value = '  +10.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.610'
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.95'
$ws.Range('E10').Value = '  +5.29%  '
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.90'
$ws.Range('E12').Value = '  +12.04%  '
$ws.Range('E13').Value = '  +1.06%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.16'
$ws.Range('E14').Value = '  +3.53%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.564.45'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '2.226.78'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.794'
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').Value = '43.595.04'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.41'
$ws.Range('E20').Value = '  -1.09%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.98'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.34'
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.04'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.70'
$ws.Range('E24').Value = '  -6.21%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.50'
$ws.Range('E26').Value = '  +10.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.14'
$ws.Range('E27').Value = '  +3.03%  '
$ws.Range('E28').Value = '  +4.94%  '
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.87'
$ws.Range('E30').Value = '  -5.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '172.58'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0914'
$ws.Range('E32').Value = '  +1.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.75'
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.111'
$ws.Range('E36').Value = '  -3.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0349'
$ws.Range('E37').Value = '  -5.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.27'
$ws.Range('E38').Value = '  -7.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.45'
$ws.Range('E39').Value = '  +14.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.56'
$ws.Range('E40').Value = '  -3.67%  '
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.214'
$ws.Range('E42').Value = '  +5.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '63.16'
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.36'
$ws.Range('E44').Value = '  -2.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.39'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0987'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '99.95'
$ws.Range('E47').Value = '  -4.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.15'
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('E49').Value = '  +2.01%  '
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('E51').Value = '  -4.72%  '
